# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
# Values that look like plain numbers are entered with a leading apostrophe so
# Excel stores them as text (matching the original inline-string cell type)
# instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.786.62"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "3.468.07"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'599.10"
$ws.Range("E5").Value = "  -3.13%  "
$ws.Range("D6").Value = "'147.40"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("D7").Value = "3.467.99"
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").Value = "'7.69"
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("D12").Value = "'0.421"
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("D14").Value = "4.053.58"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "'30.99"
$ws.Range("E15").Value = "  -6.73%  "
$ws.Range("D16").Value = "3.469.64"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "66.808.98"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "'6.36"
$ws.Range("E19").Value = "  -5.62%  "
$ws.Range("D20").Value = "'10.08"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "'15.12"
$ws.Range("E21").Value = "  -5.41%  "
$ws.Range("D22").Value = "'432.29"
$ws.Range("E22").Value = "  -4.77%  "
$ws.Range("D23").Value = "'0.603"
$ws.Range("E23").Value = "  -6.39%  "
$ws.Range("D24").Value = "'79.01"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "3.607.11"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("E27").Value = "  -9.75%  "
$ws.Range("E28").Value = "  -6.88%  "
$ws.Range("E29").Value = "  -10.18%  "
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("E31").Value = "  -7.14%  "
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'25.26"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("D35").Value = "3.457.34"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'5.88"
$ws.Range("E36").Value = "  -7.73%  "
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  -6.90%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "'7.86"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'172.77"
$ws.Range("E41").Value = "  -4.79%  "
$ws.Range("D42").Value = "'0.0879"
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("D43").Value = "'5.38"
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("E44").Value = "  -13.45%  "
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'46.32"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "'28.77"
$ws.Range("E47").Value = "  -6.76%  "
$ws.Range("E48").Value = "  -7.69%  "
$ws.Range("D49").Value = "'7.42"
$ws.Range("E49").Value = "  -4.73%  "
$ws.Range("D50").Value = "'2.37"
$ws.Range("E50").Value = "  -10.48%  "
$ws.Range("D51").Value = "'0.963"
$ws.Range("E51").Value = "  -5.38%  "
